{"js": "// Apply strikethrough formatting to the \"Further explanation required:\" section\n// (the heading plus its first bullet group, down to the paragraph ending\n// \"...biophysical Ca2+ model?\"), matching the target revision.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Identify the contiguous block of paragraphs to strike through by their\n// (trimmed) text content, so the script does not depend on fixed indices.\nconst startMarker = \"Further explanation required:\";\nconst endMarker = \"How do these parameter and conversion choices alter the conclusions drawn about the biophysical Ca2+ model?\";\n\nconst items = paragraphs.items;\nlet startIndex = -1;\nlet endIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text.trim();\n  if (startIndex === -1 && text === startMarker) {\n    startIndex = i;\n  } else if (startIndex !== -1 && text === endMarker) {\n    endIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1 || endIndex === -1) {\n  throw new Error(\"Could not locate the 'Further explanation required' block.\");\n}\n\nfor (let i = startIndex; i <= endIndex; i++) {\n  const paragraph = items[i];\n  // Skip blank separator paragraphs (e.g. the empty spacer line between the\n  // heading and the bulleted list) \u2014 those are untouched by the revision.\n  if (paragraph.text.trim().length === 0) {\n    continue;\n  }\n  // Setting strikeThrough on the paragraph's own font applies it both to the\n  // paragraph mark run properties and to every run contained in it.\n  paragraph.font.strikeThrough = true;\n}\n\nawait context.sync();\n", "ps1": "# Apply strikethrough formatting to the \"Further explanation required:\" section\n# (the heading plus its first bullet group, down to the paragraph ending\n# \"...biophysical Ca2+ model?\"), matching the target revision.\n\n$doc = $word.ActiveDocument\n\n$startMarker = \"Further explanation required:\"\n$endMarker = \"How do these parameter and conversion choices alter the conclusions drawn about the biophysical Ca2+ model?\"\n\n# Locate the contiguous block of paragraphs to strike through by their\n# (trimmed) text content, so the script does not depend on fixed indices.\n$startIndex = -1\n$endIndex = -1\n$count = $doc.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $doc.Paragraphs.Item($i).Range.Text.Trim()\n    if ($startIndex -eq -1 -and $text -eq $startMarker) {\n        $startIndex = $i\n    } elseif ($startIndex -ne -1 -and $text -eq $endMarker) {\n        $endIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1 -or $endIndex -eq -1) {\n    throw \"Could not locate the 'Further explanation required' block.\"\n}\n\nfor ($i = $startIndex; $i -le $endIndex; $i++) {\n    $p = $doc.Paragraphs.Item($i)\n    # Skip blank separator paragraphs (e.g. the empty spacer line between the\n    # heading and the bulleted list) \u2014 those are untouched by the revision.\n    if ($p.Range.Text.Trim().Length -eq 0) {\n        continue\n    }\n    # Setting StrikeThrough on the paragraph's own Range applies it both to\n    # the paragraph mark run properties and to every run contained in it.\n    $p.Range.Font.StrikeThrough = 1\n}\n"}
